$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire second row (the original first data row, dated 25569),
# shifting all subsequent rows up by one.
$ws.Rows.Item(2).Delete()

# Delete the entire "index" column (K), which is no longer needed.
$ws.Columns.Item(11).Delete()
